$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = '[Feng%Pan%NULL%0,                  Lian%Yang%NULL%2,                  Yuncheng%Li%NULL%2,                  Bo%Liang%NULL%2,                  Lin%Li%NULL%3,                  Tianhe%Ye%NULL%2,                  Lingli%Li%NULL%2,                  Dehan%Liu%NULL%2,                  Shan%Gui%NULL%2,                  Yu%Hu%NULL%14,                  Chuansheng%Zheng%NULL%0]'
$ws.Range("E3").Value = '[Barnaby%D.P.%coreGivesNoEmail%0,                 Becker%L.B.%coreGivesNoEmail%2,                 Chelico%J.D.%coreGivesNoEmail%2,                 Cohen%S.L.%coreGivesNoEmail%2,                 Cookingham%J.%coreGivesNoEmail%2,                 Coppa%K.%coreGivesNoEmail%2,                 Crawford%J.M.%coreGivesNoEmail%2,                 Davidson%K.W.%coreGivesNoEmail%2,                 Diefenbach%M.A.%coreGivesNoEmail%2,                 Dominello%A.J.%coreGivesNoEmail%2,                 Duer-Hefele%J.%coreGivesNoEmail%2,                 Falzon%L.%coreGivesNoEmail%2,                 Gitlin%J.%coreGivesNoEmail%2,                 Hajizadeh%N.%coreGivesNoEmail%2,                 Harvin%T.G.%coreGivesNoEmail%2,                 Hirsch%J.S.%coreGivesNoEmail%2,                 Hirschwerk%D.A.%coreGivesNoEmail%2,                 Kim%E.J.%coreGivesNoEmail%2,                 Kozel%Z.M.%coreGivesNoEmail%2,                 Marrast%L.M.%coreGivesNoEmail%2,                 McGinn%T.%coreGivesNoEmail%2,                 Mogavero%J.N.%coreGivesNoEmail%2,                 Narasimhan%M.%coreGivesNoEmail%2,                 Osorio%G.A.%coreGivesNoEmail%2,                 Qiu%M.%coreGivesNoEmail%2,                 Richardson%S.%coreGivesNoEmail%2,                 Zanos%T.P.%coreGivesNoEmail%2]'
$ws.Range("E4").Value = '[Mario%Rivera-Izquierdo%NULL%2,                  María%del Carmen Valero-Ubierna%NULL%4,                  María%del Carmen Valero-Ubierna%NULL%0,                  Juan Luis%R-delAmo%NULL%2,                  Miguel Ángel%Fernández-García%NULL%2,                  Silvia%Martínez-Diz%NULL%2,                  Arezu%Tahery-Mahmoud%NULL%2,                  Marta%Rodríguez-Camacho%NULL%2,                  Ana Belén%Gámiz-Molina%NULL%2,                  Nicolás%Barba-Gyengo%NULL%2,                  Pablo%Gámez-Baeza%NULL%2,                  Celia%Cabrero-Rodríguez%NULL%2,                  Pedro Antonio%Guirado-Ruiz%NULL%2,                  Divina Tatiana%Martín-Romero%NULL%2,                  Antonio Jesús%Láinez-Ramos-Bossini%NULL%2,                  María Rosa%Sánchez-Pérez%NULL%2,                  José%Mancera-Romero%NULL%2,                  Miguel%García-Martín%NULL%2,                  Luis Miguel%Martín-delosReyes%NULL%2,                  Virginia%Martínez-Ruiz%NULL%4,                  Virginia%Martínez-Ruiz%NULL%0,                  Pablo%Lardelli-Claret%NULL%2,                  Eladio%Jiménez-Mejías%NULL%2,                  Muhammad%Adrish%NULL%4,                  Muhammad%Adrish%NULL%0,                  NULL%NULL%NULL%25,                  NULL%NULL%NULL%0]'
$ws.Range("E5").Value = '[Qiurong%Ruan%NULL%0,                  Kun%Yang%NULL%0,                  Kun%Yang%NULL%0,                  Wenxia%Wang%NULL%0,                  Wenxia%Wang%NULL%0,                  Lingyu%Jiang%NULL%0,                  Lingyu%Jiang%NULL%0,                  Jianxin%Song%songsingsjx@sina.com%0,                  Jianxin%Song%songsingsjx@sina.com%0]'
$ws.Range("E6").Value = '[Grace%Salacup%salacupg@einstein.edu%3,                  Kevin Bryan%Lo%NULL%5,                  Kevin Bryan%Lo%NULL%0,                  Fahad%Gul%NULL%2,                  Eric%Peterson%NULL%3,                  Robert%De Joy%NULL%2,                  Ruchika%Bhargav%NULL%3,                  Jerald%Pelayo%NULL%3,                  Jeri%Albano%NULL%3,                  Zurab%Azmaiparashvili%NULL%3,                  Sadia%Benzaquen%NULL%3,                  Gabriel%Patarroyo‐Aponte%NULL%2,                  Janani%Rangaswami%NULL%3]'
$ws.Range("E8").Value = '[Yufeng%Shang%NULL%2,                  Tao%Liu%NULL%4,                  Yongchang%Wei%NULL%2,                  Jingfeng%Li%NULL%2,                  Liang%Shao%NULL%2,                  Minghui%Liu%NULL%2,                  Yongxi%Zhang%NULL%0,                  Zhigang%Zhao%NULL%2,                  Haibo%Xu%NULL%5,                  Zhiyong%Peng%NULL%5,                  Fuling%Zhou%zhoufuling@whu.edu.cn%2,                  Xinghuan%Wang%wangxinghuan@whu.edu.cn%4]'
$ws.Range("E9").Value = '[Shaobo%Shi%NULL%4,                  Mu%Qin%qinmuae@163.com%2,                  Yuli%Cai%NULL%2,                  Tao%Liu%NULL%0,                  Bo%Shen%NULL%2,                  Fan%Yang%NULL%0,                  Sheng%Cao%NULL%2,                  Xu%Liu%NULL%4,                  Xu%Liu%NULL%0,                  Yaozu%Xiang%NULL%2,                  Qinyan%Zhao%NULL%2,                  He%Huang%huanghe1977@whu.edu.cn%2,                  Bo%Yang%yybb112@whu.edu.cn%0,                  Congxin%Huang%NULL%2]'
$ws.Range("E10").Value = '[Rita de Cássia Menezes%Soares%NULL%2,                  Larissa Rodrigues%Mattos%NULL%2,                  Letícia Martins%Raposo%NULL%2]'
$ws.Range("E11").Value = '[Haiying%Sun%NULL%4,                  Ruoqi%Ning%NULL%3,                  Yu%Tao%NULL%4,                  Chong%Yu%NULL%3,                  Xiaoyan%Deng%NULL%3,                  Caili%Zhao%NULL%3,                  Silu%Meng%NULL%3,                  Fangxu%Tang%89650793@qq.com%3,                  Dong%Xu%89650793@qq.com%11,                  Dong%Xu%89650793@qq.com%0]'
$ws.Range("E12").Value = '[Kun%Wang%NULL%2,                  Peiyuan%Zuo%NULL%2,                  Yuwei%Liu%NULL%2,                  Meng%Zhang%NULL%2,                  Xiaofang%Zhao%NULL%2,                  Songpu%Xie%NULL%2,                  Hao%Zhang%NULL%2,                  Xinglin%Chen%NULL%2,                  Chengyun%Liu%chengyunliu@hust.edu.cn%2]'
$ws.Range("E13").Value = '[Bo%XU%NULL%2,                  Cun-yu%FAN%NULL%2,                  An-lu%WANG%NULL%2,                  Yi-long%ZOU%NULL%2,                  Yi-han%YU%NULL%2,                  Cong%HE%NULL%2,                  Wen-guang%XIA%NULL%2,                  Ji-xian%ZHANG%NULL%2,                  Qing%MIAO%NULL%2]'
$ws.Range("E14").Value = '[Xisheng%Yan%NULL%3,                  Fen%Li%NULL%2,                  Xiao%Wang%NULL%3,                  Jie%Yan%NULL%3,                  Fen%Zhu%NULL%2,                  Shifan%Tang%NULL%2,                  Yingzhong%Deng%NULL%2,                  Hua%Wang%NULL%3,                  Rui%Chen%NULL%2,                  Zhili%Yu%NULL%2,                  Yaping%Li%NULL%2,                  Jingzhou%Shang%NULL%2,                  Lingjun%Zeng%NULL%2,                  Jie%Zhao%NULL%2,                  Chaokun%Guan%NULL%2,                  Qiaomei%Liu%NULL%2,                  Haifeng%Chen%NULL%2,                  Wei%Gong%NULL%2,                  Xin%Huang%NULL%2,                  Yu‐Jiao%Zhang%NULL%2,                  Jianguang%Liu%NULL%2,                  Xiaoyan%Dong%NULL%2,                  Wen%Zheng%zhengwen12@mails.jlu.edu.cn%2,                  Shaoping%Nie%spnie@126.com%2,                  Dongsheng%Li%dongshengli196809@163.com%2]'
$ws.Range("E15").Value = '[Qing%Yang%NULL%2,                  Ying%Zhou%NULL%2,                  Xinrong%Wang%NULL%2,                  Shan%Gao%NULL%2,                  Yang%Xiao%NULL%2,                  Weiming%Zhang%NULL%2,                  Yi%Hu%huyizxyy@163.com%0,                  Yafei%Wang%wyf_527@163.com%2]'
$ws.Range("E16").Value = '[Xiaobo%Yang%NULL%0,                  Qingyu%Yang%NULL%2,                  Yaxin%Wang%NULL%0,                  Yongran%Wu%NULL%0,                  Jiqian%Xu%NULL%0,                  Yuan%Yu%NULL%0,                  You%Shang%NULL%0]'
$ws.Range("E17").Value = '[Wenjing%Ye%yewenjing@xinhuamed.com.cn%2,                  Guoxi%Chen%4799082@qq.com%2,                  Xiaopan%Li%xiaopanli0224@126.com%2,                  Xing%Lan%474137452@qq.com%2,                  Chen%Ji%c.ji.3@warwick.ac.uk%2,                  Min%Hou%min-hou710@hotmail.com%2,                  Di%Zhang%zhangdizhangdi1234@163.com%2,                  Guangwang%Zeng%785663609@qq.com%2,                  Yaling%Wang%953822402@qq.com%2,                  Cheng%Xu%xucheng112358@126.com%2,                  Weiwei%Lu%luweiwei100@126.com%2,                  Ruolin%Cui%crlqwerty@163.com%2,                  Yuyang%Cai%caiyuyang@sjtu.edu.cn%2,                  Hai%Huang%1220775601@qq.com%2,                  Ling%Yang%yangling01@xinhuamed.com.cn%3]'
$ws.Range("E18").Value = '[Caizheng%Yu%NULL%2,                  Qing%Lei%NULL%2,                  Wenkai%Li%NULL%2,                  Xiong%Wang%NULL%0,                  Wei%Liu%NULL%0,                  Xionglin%Fan%NULL%2,                  Wengang%Li%228907211@qq.com%2]'
$ws.Range("E19").Value = '[Fei%Zhou%NULL%0,                  Ting%Yu%NULL%0,                  Ronghui%Du%NULL%0,                  Guohui%Fan%NULL%0,                  Ying%Liu%NULL%0,                  Zhibo%Liu%NULL%0,                  Jie%Xiang%NULL%0,                  Yeming%Wang%NULL%0,                  Bin%Song%NULL%0,                  Xiaoying%Gu%NULL%0,                  Lulu%Guan%NULL%0,                  Yuan%Wei%NULL%0,                  Hui%Li%NULL%0,                  Xudong%Wu%NULL%0,                  Jiuyang%Xu%NULL%0,                  Shengjin%Tu%NULL%0,                  Yi%Zhang%NULL%0,                  Hua%Chen%NULL%0,                  Bin%Cao%NULL%0]'
$ws.Range("E20").Value = '[Moran%Amit%NULL%2,                  Alex%Sorkin%NULL%2,                  Jacob%Chen%NULL%2,                  Barak%Cohen%NULL%4,                  Barak%Cohen%NULL%0,                  Dana%Karol%NULL%4,                  Dana%Karol%NULL%0,                  Avishai M%Tsur%NULL%2,                  Shaul%Lev%NULL%4,                  Shaul%Lev%NULL%0,                  Tal%Rozenblat%NULL%2,                  Ayana%Dvir%NULL%4,                  Ayana%Dvir%NULL%0,                  Geva%Landau%NULL%2,                  Lidar%Fridrich%NULL%4,                  Lidar%Fridrich%NULL%0,                  Elon%Glassberg%NULL%2,                  Shani%Kesari%NULL%2,                  Sigal%Sviri%NULL%2,                  Ram%Gelman%NULL%2,                  Asaf%Miller%NULL%2,                  Danny%Epstein%NULL%2,                  Ronny%Ben-Avi%NULL%4,                  Ronny%Ben-Avi%NULL%0,                  Moshe%Matan%NULL%2,                  Daniel J.%Jakobson%NULL%4,                  Daniel J.%Jakobson%NULL%0,                  Tarif%Bader%NULL%4,                  Tarif%Bader%NULL%0,                  David%Dahan%NULL%4,                  David%Dahan%NULL%0,                  Daniel A.%King%NULL%2,                  Anat%Ben-Ari%NULL%2,                  Arie%Soroksky%NULL%2,                  Alon%Bar%NULL%4,                  Alon%Bar%NULL%0,                  Noam%Fink%NULL%2,                  Pierre%Singer%NULL%2,                  Avi%Benov%NULL%2]'
$ws.Range("E21").Value = '[Alexander%Muacevic%NULL%0,                  John R%Adler%NULL%2,                  Muhammad Sohaib%Asghar%NULL%4,                  Muhammad Sohaib%Asghar%NULL%0,                  Syed Jawad%Haider Kazmi%NULL%2,                  Noman%Ahmed Khan%NULL%2,                  Mohammed%Akram%NULL%2,                  Salman%Ahmed Khan%NULL%2,                  Uzma%Rasheed%NULL%2,                  Maira%Hassan%NULL%2,                  Gul Muhammad%Memon%NULL%2]'
$ws.Range("E22").Value = '[Pedro%Baqui%NULL%2,                  Ioana%Bica%NULL%2,                  Valerio%Marra%marra@cosmo-ufes.org%2,                  Ari%Ercole%NULL%2,                  Mihaela%van der Schaar%NULL%2]'
$ws.Range("E23").Value = '[Anesi%Adriano%coreGivesNoEmail%3,                 Bettinardi%Alessandra%coreGivesNoEmail%3,                 Bonetti%Graziella%coreGivesNoEmail%3,                 Borrelli%Gianluca%coreGivesNoEmail%3,                 Fiordalisi%Gianfranco%coreGivesNoEmail%3,                 Lippi%Giuseppe%coreGivesNoEmail%3,                 Manelli%Filippo%coreGivesNoEmail%3,                 Marino%Antonio%coreGivesNoEmail%3,                 Menolfi%Annamaria%coreGivesNoEmail%3,                 Patroni%Andrea%coreGivesNoEmail%3,                 Saggini%Sara%coreGivesNoEmail%3,                 Volpi%Roberta%coreGivesNoEmail%3]'
$ws.Range("E24").Value = '[Andrea%Borghesi%NULL%2,                  Angelo%Zigliani%NULL%2,                  Salvatore%Golemi%NULL%2,                  Nicola%Carapella%NULL%2,                  Patrizia%Maculotti%NULL%2,                  Davide%Farina%NULL%2,                  Roberto%Maroldi%NULL%2]'
$ws.Range("E25").Value = '[Alberto M.%Borobia%NULL%2,                  Antonio J.%Carcas%NULL%4,                  Antonio J.%Carcas%NULL%0,                  Francisco%Arnalich%NULL%2,                  Rodolfo%Álvarez-Sala%NULL%4,                  Rodolfo%Álvarez-Sala%NULL%0,                  Jaime%Monserrat-Villatoro%NULL%2,                  Manuel%Quintana%NULL%2,                  Juan Carlos%Figueira%NULL%2,                  Rosario M.%Torres Santos-Olmo%NULL%2,                  Julio%García-Rodríguez%NULL%4,                  Julio%García-Rodríguez%NULL%0,                  Alberto%Martín-Vega%NULL%2,                  Antonio%Buño%NULL%2,                  Elena%Ramírez%NULL%2,                  Gonzalo%Martínez-Alés%NULL%4,                  Gonzalo%Martínez-Alés%NULL%0,                  Nicolás%García-Arenzana%NULL%2,                  M. Concepción%Núñez%NULL%4,                  M. Concepción%Núñez%NULL%0,                  Milagros%Martí-de-Gracia%NULL%2,                  Francisco%Moreno Ramos%NULL%2,                  Francisco%Reinoso-Barbero%NULL%2,                  Alejandro%Martin-Quiros%NULL%2,                  Angélica%Rivera Núñez%NULL%2,                  Jesús%Mingorance%NULL%2,                  Carlos J.%Carpio Segura%NULL%4,                  Carlos J.%Carpio Segura%NULL%0,                  Daniel%Prieto Arribas%NULL%2,                  Esther%Rey Cuevas%NULL%2,                  Concepción%Prados Sánchez%NULL%2,                  Juan J.%Rios%NULL%2,                  Miguel A.%Hernán%NULL%2,                  Jesús%Frías%NULL%2,                  José R.%Arribas%NULL%2,                  NULL%NULL%NULL%0]'
$ws.Range("E26").Value = '[Simon E.%Brill%simon.brill@nhs.net%0,                  Hannah C.%Jarvis%NULL%4,                  Hannah C.%Jarvis%NULL%0,                  Ezgi%Ozcan%NULL%2,                  Thomas L. P.%Burns%NULL%2,                  Rabia A.%Warraich%NULL%2,                  Lisa J.%Amani%NULL%2,                  Amina%Jaffer%NULL%2,                  Stephanie%Paget%NULL%2,                  Anand%Sivaramakrishnan%NULL%2,                  Dean D.%Creer%NULL%2]'
$ws.Range("E27").Value = '[Jianlei%Cao%NULL%0,                  Wen-Jun%Tu%tuwenjun@irm-cams.ac.cn%0,                  Wenlin%Cheng%NULL%0,                  Lei%Yu%NULL%0,                  Ya-Kun%Liu%NULL%5,                  Xiaoyong%Hu%NULL%5,                  Qiang%Liu%NULL%0]'
$ws.Range("E28").Value = '[B.%Carter%NULL%2,                  J.T.%Collins%NULL%2,                  F.%Barlow-Pay%NULL%2,                  F.%Rickard%NULL%2,                  E.%Bruce%NULL%2,                  A.%Verduri%NULL%2,                  T.J.%Quinn%NULL%2,                  E.%Mitchell%NULL%2,                  A.%Price%NULL%2,                  A.%Vilches-Moraga%NULL%2,                  M.J.%Stechman%NULL%2,                  R.%Short%NULL%2,                  A.%Einarsson%NULL%2,                  P.%Braude%NULL%2,                  S.%Moug%NULL%2,                  P.K.%Myint%NULL%2,                  J.%Hewitt%NULL%2,                  L.%Pearce%NULL%2,                  K.%McCarthy%NULL%2,                  C.%Davey%NULL%2,                  S.%Jones%NULL%2,                  K.%Lunstone%NULL%2,                  A.%Cavenagh%NULL%2,                  C.%Silver%NULL%2,                  T.%Telford%NULL%2,                  R.%Simmons%NULL%2,                  M.%Holloway%NULL%2,                  J.%Hesford%NULL%2,                  T.%El Jichi Mutasem%NULL%2,                  S.%Singh%NULL%2,                  D.%Paxton%NULL%2,                  W.%Harris%NULL%2,                  N.%Galbraith%NULL%2,                  E.%Bhatti%NULL%2,                  J.%Edwards%NULL%2,                  S.%Duffy%NULL%2,                  J.%Kelly%NULL%2,                  C.%Murphy%NULL%2,                  C.%Bisset%NULL%2,                  R.%Alexander%NULL%2,                  M.%Garcia%NULL%2,                  S.%Sangani%NULL%2,                  T.%Kneen%NULL%2,                  T.%Lee%NULL%2,                  A.%McGovern%NULL%2,                  G.%Guaraldi%NULL%2,                  E.%Clini%NULL%2]'
$ws.Range("E29").Value = '[Fuyang%Chen%NULL%3,                  Wenwu%Sun%NULL%2,                  Shengrong%Sun%NULL%2,                  Zhiyu%Li%lizhiyu@whu.edu.cn%2,                  Zhong%Wang%zhongwangchn@whu.edu.cn%2,                  Li%Yu%yuliwhzxyy@163.com%4,                  Li%Yu%yuliwhzxyy@163.com%0]'
$ws.Range("E30").Value = '[Ruchong%Chen%NULL%0,                  Ling%Sang%NULL%7,                  Mei%Jiang%NULL%4,                  Zhaowei%Yang%NULL%2,                  Nan%Jia%NULL%2,                  Wanyi%Fu%NULL%2,                  Jiaxing%Xie%NULL%2,                  Weijie%Guan%NULL%4,                  Wenhua%Liang%NULL%4,                  Zhengyi%Ni%NULL%4,                  Yu%Hu%NULL%0,                  Lei%Liu%NULL%0,                  Hong%Shan%NULL%0,                  Chunliang%Lei%NULL%4,                  Yixiang%Peng%NULL%4,                  Li%Wei%NULL%12,                  Yong%Liu%NULL%12,                  Yahua%Hu%NULL%4,                  Peng%Peng%NULL%16,                  Jianming%Wang%NULL%0,                  Jiyang%Liu%NULL%4,                  Zhong%Chen%NULL%12,                  Gang%Li%NULL%13,                  Zhijian%Zheng%NULL%4,                  Shaoqin%Qiu%NULL%4,                  Jie%Luo%NULL%12,                  Changjiang%Ye%NULL%4,                  Shaoyong%Zhu%NULL%4,                  Jinping%Zheng%NULL%4,                  Nuofu%Zhang%NULL%4,                  Yimin%Li%NULL%5,                  Jianxing%He%NULL%4,                  Jing%Li%NULL%3,                  Shiyue%Li%NULL%4,                  Nanshan%Zhong%NULL%5,                  NULL%NULL%NULL%0]'
$ws.Range("E31").Value = '[Tao%Chen%NULL%0,                  Di%Wu%NULL%5,                  Huilong%Chen%NULL%5,                  Weiming%Yan%NULL%5,                  Danlei%Yang%NULL%5,                  Guang%Chen%NULL%5,                  Ke%Ma%NULL%6,                  Dong%Xu%NULL%0,                  Haijing%Yu%NULL%5,                  Hongwu%Wang%NULL%5,                  Tao%Wang%NULL%18,                  Wei%Guo%NULL%5,                  Jia%Chen%NULL%5,                  Chen%Ding%NULL%5,                  Xiaoping%Zhang%NULL%5,                  Jiaquan%Huang%NULL%5,                  Meifang%Han%NULL%5,                  Shusheng%Li%NULL%6,                  Xiaoping%Luo%NULL%5,                  Jianping%Zhao%NULL%0,                  Qin%Ning%NULL%5]'
$ws.Range("E32").Value = '[Anying%Cheng%NULL%0,                  Liu%Hu%NULL%2,                  Yiru%Wang%NULL%2,                  Luyan%Huang%NULL%2,                  Lingxi%Zhao%NULL%2,                  Congcong%Zhang%NULL%2,                  Xiyue%Liu%NULL%2,                  Ranran%Xu%NULL%2,                  Feng%Liu%NULL%2,                  Jinping%Li%NULL%2,                  Dawei%Ye%NULL%2,                  Tao%Wang%NULL%0,                  Yongman%Lv%lvyongman@126.com%2,                  Qingquan%Liu%qqliutj@163.com%2]'
$ws.Range("E33").Value = '[Fabio%Ciceri%ciceri.fabio@hsr.it%2,                  Antonella%Castagna%NULL%2,                  Patrizia%Rovere-Querini%NULL%3,                  Francesco%De Cobelli%NULL%2,                  Annalisa%Ruggeri%NULL%3,                  Laura%Galli%NULL%2,                  Caterina%Conte%NULL%2,                  Rebecca%De Lorenzo%NULL%2,                  Andrea%Poli%NULL%2,                  Alberto%Ambrosio%NULL%2,                  Carlo%Signorelli%NULL%3,                  Eleonora%Bossi%NULL%2,                  Maria%Fazio%NULL%2,                  Cristina%Tresoldi%NULL%2,                  Sergio%Colombo%NULL%2,                  Giacomo%Monti%NULL%2,                  Efgeny%Fominskiy%NULL%2,                  Stefano%Franchini%NULL%2,                  Marzia%Spessot%NULL%2,                  Carlo%Martinenghi%NULL%2,                  Michele%Carlucci%NULL%2,                  Luigi%Beretta%NULL%2,                  Anna Maria%Scandroglio%NULL%2,                  Massimo%Clementi%NULL%2,                  Massimo%Locatelli%NULL%2,                  Moreno%Tresoldi%NULL%3,                  Paolo%Scarpellini%NULL%3,                  Gianvito%Martino%NULL%2,                  Emanuele%Bosi%NULL%2,                  Lorenzo%Dagna%NULL%3,                  Adriano%Lazzarin%NULL%2,                  Giovanni%Landoni%NULL%3,                  Alberto%Zangrillo%NULL%3]'
$ws.Range("E34").Value = '[Yan%Deng%NULL%0,                  Wei%Liu%NULL%0,                  Kui%Liu%NULL%0,                  Yuan-Yuan%Fang%NULL%0,                  Jin%Shang%NULL%2,                  Ling%Zhou%NULL%0,                  Ke%Wang%NULL%3,                  Fan%Leng%NULL%2,                  Shuang%Wei%NULL%0,                  Lei%Chen%NULL%2,                  Hui-Guo%Liu%NULL%0,                  Pei-Fang%Wei%NULL%0,                  Pei-Fang%Wei%NULL%0]'
$ws.Range("E35").Value = '[Rong-Hui%Du%NULL%0,                  Li-Rong%Liang%NULL%2,                  Cheng-Qing%Yang%NULL%2,                  Wen%Wang%NULL%4,                  Tan-Ze%Cao%NULL%2,                  Ming%Li%NULL%2,                  Guang-Yun%Guo%NULL%2,                  Juan%Du%NULL%2,                  Chun-Lan%Zheng%NULL%2,                  Qi%Zhu%NULL%2,                  Ming%Hu%NULL%2,                  Xu-Yan%Li%NULL%4,                  Peng%Peng%NULL%0,                  Huan-Zhong%Shi%NULL%4]'
$ws.Range("E36").Value = '[Shan%Gao%NULL%2,                  Fang%Jiang%NULL%4,                  Wei%Jin%NULL%2,                  Yuan%Shi%NULL%3,                  Leilei%Yang%NULL%2,                  Yanqiong%Xia%NULL%2,                  Linyan%Jia%NULL%2,                  Bo%Wang%NULL%2,                  Han%Lin%NULL%2,                  Yin%Cai%NULL%2,                  Zhengyuan%Xia%NULL%4,                  Jian%Peng%NULL%2]'
$ws.Range("E37").Value = '[Pedro David%Wendel Garcia%pedrodavid.wendelgarcia@usz.ch%2,                  Thierry%Fumeaux%thierry.fumeaux@ghol.ch%2,                  Philippe%Guerci%NULL%2,                  Dorothea Monika%Heuberger%dorotheamonika.heuberger@usz.ch%2,                  Jonathan%Montomoli%NULL%2,                  Ferran%Roche-Campo%NULL%2,                  Reto Andreas%Schuepbach%reto.schuepbach@usz.ch%2,                  Matthias Peter%Hilty%matthias.hilty@usz.ch%2,                  Mario%Alfaro Farias%NULL%2,                  Antoni%Margarit%NULL%2,                  Gerardo%Vizmanos-Lamotte%NULL%2,                  Thomas%Tschoellitsch%NULL%2,                  Jens%Meier%NULL%2,                  Francesco S.%Cardona%NULL%2,                  Josef%Skola%NULL%2,                  Lenka%Horakova%NULL%2,                  Hernan%Aguirre-Bermeo%NULL%2,                  Janina%Apolo%NULL%2,                  Emmanuel%Novy%NULL%2,                  Marie-Reine%Losser%NULL%2,                  Geoffrey%Jurkolow%NULL%2,                  Gauthier%Delahaye%NULL%2,                  Sascha%David%NULL%2,                  Tobias%Welte%NULL%2,                  Tobias%Wengenmayer%NULL%2,                  Dawid L.%Staudacher%NULL%2,                  Theodoros%Aslanidis%NULL%2,                  Barna%Babik%NULL%2,                  Anita%Korsos%NULL%2,                  Janos%Gal%NULL%2,                  Hermann%Csaba%NULL%2,                  Abele%Donati%NULL%2,                  Andrea%Carsetti%NULL%2,                  Fabrizio%Turrini%NULL%2,                  Maria Sole%Simonini%NULL%2,                  Roberto%Ceriani%NULL%2,                  Martina%Murrone%NULL%2,                  Emanuele%Rezoagli%NULL%2,                  Giovanni%Vitale%NULL%2,                  Alberto%Fogagnolo%NULL%2,                  Savino%Spadaro%NULL%2,                  Maddalena Alessandra%Wu%NULL%2,                  Chiara%Cogliati%NULL%2,                  Riccardo%Colombo%NULL%2,                  Emanuele%Catena%NULL%2,                  Francesca%Facondini%NULL%2,                  Antonella%Potalivo%NULL%2,                  Gianfilippo%Gangitano%NULL%2,                  Tiziana%Perin%NULL%2,                  Maria Grazia%Bocci%NULL%2,                  Massimo%Antonelli%NULL%2,                  Diederik%Gommers%NULL%2,                  Can%Ince%NULL%2,                  Eric%Mayor-Vázquez%NULL%2,                  Maria%Cruz%NULL%2,                  Martin%Delgado%NULL%2,                  Raquel Rodriguez%Garcia%NULL%2,                  Jorge%Gamez Zapata%NULL%2,                  Begoña%Zalba-Etayo%NULL%2,                  Herminia%Lozano-Gomez%NULL%2,                  Pedro%Castro%NULL%2,                  Adrian%Tellez%NULL%2,                  Adriana%Jacas%NULL%2,                  Guido%Muñoz%NULL%2,                  Rut%Andrea%NULL%2,                  Jose%Ortiz%NULL%2,                  Eduard%Quintana%NULL%2,                  Irene%Rovira%NULL%2,                  Enric%Reverter%NULL%2,                  Javier%Fernandez%NULL%2,                  Miquel%Ferrer%NULL%2,                  Joan R.%Badia%NULL%2,                  Arantxa%Lander Azcona%NULL%2,                  Jesus Escos%Orta%NULL%2,                  Philipp%Bühler%NULL%2,                  Silvio%Brugger%NULL%2,                  Daniel%Hofmaenner%NULL%2,                  Simone%Unseld%NULL%2,                  Frank%Ruschitzka%NULL%2,                  Mallory%Moret-Bochatay%NULL%2,                  Bernd%Yuen%NULL%2,                  Thomas%Hillermann%NULL%2,                  Hatem%Ksouri%NULL%2,                  Govind Oliver%Sridharan%NULL%2,                  Anette%Ristic%NULL%2,                  Michael%Sepulcri%NULL%2,                  Miodrag%Filipovic%NULL%2,                  Urs%Pietsch%NULL%2,                  Petra%Salomon%NULL%2,                  Iris%Drvaric%NULL%2,                  Peter%Schott%NULL%2,                  Severin%Urech%NULL%2,                  Adriana%Lambert%NULL%2,                  Lukas%Merki%NULL%2,                  Marcus%Laube%NULL%2,                  Frank%Hillgaertner%NULL%2,                  Marianne%Sieber%NULL%2,                  Alexander%Dullenkopf%NULL%2,                  Lina%Petersen%NULL%2,                  Serge%Grazioli%NULL%2,                  Peter C.%Rimensberger%NULL%2,                  Isabelle%Fleisch%NULL%2,                  Jerome%Lavanchy%NULL%2,                  Katharina%Marquardt%NULL%2,                  Karim%Shaikh%NULL%2,                  Hermann%Redecker%NULL%2,                  Michael%Stephan%NULL%2,                  Jan%Brem%NULL%2,                  Bjarte%Rogdo%NULL%2,                  Andre%Birkenmaier%NULL%2,                  Friederike%Meyer zu Bentrup%NULL%2,                  Patricia%Fodor%NULL%2,                  Pascal%Locher%NULL%2,                  Giovanni%Camen%NULL%2,                  Martin%Siegemund%NULL%2,                  Nuria%Zellweger%NULL%2,                  Marie-Madlen%Jeitziner%NULL%2,                  Beatrice%Jenni-Moser%NULL%2,                  Christian%Bürkle%NULL%2,                  Gian-Reto%Kleger%NULL%2,                  Marilene%Franchitti Laurent%NULL%2,                  Jean-Christophe%Laurent%NULL%2,                  Tomislav%Gaspert%NULL%2,                  Marija%Jovic%NULL%2,                  Michael%Studhalter%NULL%2,                  Christoph%Haberthuer%NULL%2,                  Roger F.%Lussman%NULL%2,                  Daniela%Selz%NULL%2,                  Didier%Naon%NULL%2,                  Romano%Mauri%NULL%2,                  Samuele%Ceruti%NULL%2,                  Julien%Marrel%NULL%2,                  Mirko%Brenni%NULL%2,                  Rolf%Ensner%NULL%2,                  Nadine%Gehring%NULL%2,                  Antje%Heise%NULL%2,                  Tobias%Huebner%NULL%2,                  Thomas A.%Neff%NULL%2,                  Sara%Cereghetti%NULL%2,                  Filippo%Boroli%NULL%2,                  Jerome%Pugin%NULL%2,                  Nandor%Marczin%NULL%2,                  Joyce%Wong%NULL%2]'
$ws.Range("E38").Value = '[Warren%Gavin%NULL%2,                  Elliott%Campbell%NULL%2,                  Syed-Adeel%Zaidi%NULL%2,                  Neha%Gavin%NULL%2,                  Lana%Dbeibo%NULL%2,                  Cole%Beeler%NULL%2,                  Kari%Kuebler%NULL%2,                  Ahmed%Abdel-Rahman%NULL%2,                  Mark%Luetkemeyer%NULL%2,                  Areeba%Kara%NULL%2]'
$ws.Range("E39").Value = '[Vijay%Gayam%vgayam@interfaithmedical.com%2,                  Muchi Ditah%Chobufo%NULL%4,                  Muchi Ditah%Chobufo%NULL%0,                  Mohamed A.%Merghani%NULL%4,                  Mohamed A.%Merghani%NULL%0,                  Shristi%Lamichhane%NULL%2,                  Pavani Reddy%Garlapati%NULL%2,                  Mark K.%Adler%NULL%2]'
$ws.Range("E40").Value = '[Cao%Y%coreGivesNoEmail%0,                 Imam%Z%coreGivesNoEmail%2,                 Lippi%G%coreGivesNoEmail%2,                 Oran%DP%coreGivesNoEmail%2,                 Shi%S%coreGivesNoEmail%2]'
$ws.Range("E41").Value = '[Hai%Hu%huhai@wchscu.cn%2,                  Ni%Yao%NULL%4,                  Ni%Yao%NULL%0,                  Yanru%Qiu%NULL%2,                  John H.%Burton%NULL%4,                  John H.%Burton%NULL%0]'
$ws.Range("E42").Value = '[Jiaofeng%Huang%NULL%2,                  Aiguo%Cheng%NULL%2,                  Rahul%Kumar%NULL%2,                  Yingying%Fang%NULL%4,                  Yingying%Fang%NULL%0,                  Gongping%Chen%NULL%2,                  Yueyong%Zhu%NULL%2,                  Su%Lin%sumer5129@fjmu.edu.cn%4,                  Su%Lin%sumer5129@fjmu.edu.cn%0]'
$ws.Range("E43").Value = '[Jong-moon%Hwang%NULL%0,                  Ju-Hyun%Kim%NULL%2,                  Jin-Sung%Park%NULL%2,                  Min Cheol%Chang%wheel633@ynu.ac.kr%2,                  Donghwi%Park%bdome@hanmail.net%2]'
$ws.Range("E44").Value = '[K.%Khalil%NULL%2,                  K.%Agbontaen%NULL%2,                  D.%McNally%NULL%2,                  A.%Love%NULL%2,                  S.%Mandalia%NULL%2,                  W.%Banya%NULL%2,                  E.%Starren%NULL%2,                  R.%Dhunnookchand%NULL%2,                  H.%Farne%NULL%2,                  R.%Morton%NULL%2,                  G.%Davies%NULL%2,                  O.%Orhan%NULL%2,                  D%Lai%NULL%2,                  M.%Nelson%NULL%2,                  P.L.%Shah%NULL%2,                  J.L.%Garner%Justin.garner@chelwest.nhs.uk%2]'
$ws.Range("E45").Value = '[Eyal%Klang%NULL%2,                  Gassan%Kassim%NULL%2,                  Shelly%Soffer%soffer.shelly@gmail.com%2,                  Robert%Freeman%NULL%4,                  Robert%Freeman%NULL%0,                  Matthew A.%Levin%NULL%4,                  Matthew A.%Levin%NULL%0,                  David L.%Reich%NULL%2]'
$ws.Range("E46").Value = '[Sandeep%Krishnan%NULL%2,                  Kinjal%Patel%NULL%2,                  Ronak%Desai%NULL%2,                  Anupam%Sule%NULL%2,                  Peter%Paik%NULL%2,                  Ashley%Miller%NULL%2,                  Alicia%Barclay%NULL%2,                  Adam%Cassella%NULL%2,                  Jon%Lucaj%NULL%2,                  Yvonne%Royster%NULL%2,                  Joffer%Hakim%NULL%2,                  Zulfiqar%Ahmed%NULL%2,                  Farhad%Ghoddoussi%NULL%2]'
$ws.Range("E48").Value = '[ Manisha%Bhutani%null%2,        David M.%Foureau%null%2,        Shebli%Atrash%null%2,        Peter M.%Voorhees%null%2,        Saad Z.%Usmani%null%2]'
$ws.Range("E49").Value = '[Wil%Lieberman-Cribbin%NULL%2,                  Joseph%Rapp%NULL%2,                  Naomi%Alpert%NULL%2,                  Stephanie%Tuminello%NULL%2,                  Emanuela%Taioli%NULL%2]'
$ws.Range("E50").Value = '[Q.%Liu%NULL%2,                  N. C.%Song%NULL%2,                  Z. K.%Zheng%NULL%2,                  J. S.%Li%NULL%2,                  S. K.%Li%NULL%2]'
$ws.Range("E51").Value = '[Hui%Long%NULL%2,                  Lan%Nie%NULL%2,                  Xiaochen%Xiang%NULL%2,                  Huan%Li%NULL%0,                  Xiaoli%Zhang%NULL%2,                  Xiaozhi%Fu%NULL%2,                  Hongwei%Ren%NULL%2,                  Wanxin%Liu%NULL%2,                  Qiang%Wang%wangqiang@wust.edu.cn%5,                  Qingming%Wu%wuhe9224@sina.com%4,                  Qingming%Wu%wuhe9224@sina.com%0]'
$ws.Range("E52").Value = '[Miao%Luo%xref no email%3,          Jing%Liu%xref no email%0,          Weiling%Jiang%xref no email%3,          Shuang%Yue%xref no email%3,          Huiguo%Liu%xref no email%3,          Shuang%Wei%xref no email%3]'
$ws.Range("E53").Value = '[Xiaomin%Luo%luoxiaomin04@163.com%2,                  Wei%Zhou%NULL%0,                  Xiaojie%Yan%NULL%2,                  Tangxi%Guo%NULL%2,                  Benchao%Wang%NULL%2,                  Hongxia%Xia%NULL%2,                  Lu%Ye%NULL%2,                  Jun%Xiong%NULL%2,                  Zongping%Jiang%NULL%2,                  Yu%Liu%NULL%2,                  Bicheng%Zhang%NULL%2,                  Weize%Yang%NULL%2]'
$ws.Range("E54").Value = '[Ying%Luo%NULL%4,                  Liyan%Mao%NULL%4,                  Xu%Yuan%NULL%4,                  Ying%Xue%NULL%4,                  Qun%Lin%NULL%4,                  Guoxing%Tang%NULL%4,                  Huijuan%Song%NULL%4,                  Feng%Wang%fengwang@tjh.tjmu.edu.cn%4,                  Ziyong%Sun%zysun@tjh.tjmu.edu.cn%0]'
$ws.Range("E55").Value = '[Ying%Luo%NULL%0,                  Ying%Xue%NULL%0,                  Liyan%Mao%NULL%0,                  Xu%Yuan%NULL%0,                  Qun%Lin%NULL%0,                  Guoxing%Tang%NULL%0,                  Huijuan%Song%NULL%0,                  Feng%Wang%NULL%0,                  Ziyong%Sun%NULL%0]'
$ws.Range("E56").Value = '[Chiara%Masetti%NULL%2,                  Elena%Generali%NULL%2,                  Francesca%Colapietro%NULL%2,                  Antonio%Voza%NULL%2,                  Maurizio%Cecconi%NULL%2,                  Antonio%Messina%NULL%2,                  Paolo%Omodei%NULL%2,                  Claudio%Angelini%NULL%2,                  Michele%Ciccarelli%NULL%4,                  Salvatore%Badalamenti%NULL%2,                  G. Walter%Canonica%NULL%2,                  Ana%Lleo%ana.lleo@humanitas.it%2,                  Alessio%Aghemo%NULL%4,                  Alessio%Aghemo%NULL%0,                  NULL%NULL%NULL%0]'
$ws.Range("E57").Value = '[Takahisa%Mikami%NULL%2,                  Hirotaka%Miyashita%NULL%2,                  Takayuki%Yamada%NULL%2,                  Matthew%Harrington%NULL%2,                  Daniel%Steinberg%NULL%2,                  Andrew%Dunn%NULL%2,                  Evan%Siau%Evan.Siau@mountsinai.org%2]'
$ws.Range("E58").Value = '[Alexis K.%Okoh%alexis.okoh@rwjbh.org%2,                  Christoph%Sossou%NULL%4,                  Christoph%Sossou%NULL%0,                  Neha S.%Dangayach%NULL%2,                  Sherin%Meledathu%NULL%2,                  Oluwakemi%Phillips%NULL%2,                  Corinne%Raczek%NULL%2,                  Michael%Patti%NULL%2,                  Nathan%Kang%NULL%2,                  Sameer A.%Hirji%NULL%2,                  Charles%Cathcart%NULL%2,                  Christian%Engell%NULL%2,                  Marc%Cohen%NULL%2,                  Sandhya%Nagarakanti%NULL%2,                  Eliahu%Bishburg%NULL%2,                  Harpreet S.%Grewal%NULL%2]'
